$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.365.13"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.843.41"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'240.16"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'0.6288"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.07467"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'25.03"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.2898"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").Value = "'0.07729"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.845.01"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'4.980"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'0.6771"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'0.00001034"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "'81.88"
$ws.Range("D17").Value = "'6.233"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "29.388.37"
$ws.Range("D19").Value = "'229.31"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'7.402"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'158.24"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'8.516"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").Value = "'0.1354"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'0.06736"
$ws.Range("E28").Value = "  +11.45%  "
$ws.Range("D29").Value = "'1.450"
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("D30").Value = "'1.490"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "'4.063"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").Value = "'4.064"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "'1.834"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'0.6999"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'2.583"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "'2.823"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "1.237.26"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'6.777"
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("D41").Value = "'0.9414"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("D42").Value = "'0.9992"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "2.000.24"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'101.14"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "'65.54"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").Value = "'0.00000000120"
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("D47").Value = "'7.055"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'1.712"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("D49").Value = "'8.973"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "'0.1145"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "'0.3917"
$ws.Range("E51").Value = "  -0.48%  "
